$d = $word.ActiveDocument
$d.Content.Find.Execute("{{ current_date }} ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " {{ current_date }} ", 2)
